$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns: Wins (AD), Losses (AE), Ties (AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (A1) so new headers match
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team record values for all data rows (2-50)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 81  # AD = column 30
    $ws.Cells.Item($r, 31).Value = 81  # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0   # AF = column 32
}
